# Simulated Wild Card round and logged it
# Appends the new game log entries to the rolling per-game sequences (YDS, ST)
# and rolls the season totals (OFF, DEF, ST, TURNS, PEN) forward by the new game.

$wb = $excel.ActiveWorkbook

# --- YDS: append the Wild Card game log to the running OFF/DEF R(ushing) and P(assing) yardage sequences ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value() + ' 2 40 0 0 -2 5 5 6 0 3 1 1 6 0 8 15 7 6 2 11 2 32 2 2 1 1 12 6 -2 -2 4 26 15 1 7 6 5 9 3 8 2 4 9 1 3 9 16 2 15 10 1 3 7 6 4 0'
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value() + ' -1 10 0 4 2 0 1 4 -1 20 2 5 -2 1 3 1 -1 2 1 16 5 2 2 3 2 3 14 14 0 0 1 8 2 7 1 5 1'
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value() + ' 17 6 10 10 19 11 6 2 5 3 9 13 17 3 12 14 8 3 28 9 4 27 -2 5 15 9 0 7 8 24 2 19 6 11 22 19 45 4 9 4 34 8 19 19 38 1'
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value() + ' 16 9 5 5 2 10 9 30 4 12 10 2 5 8 7 10 7 1 43 4 7 15 3 6 2 18'

# --- OFF: season totals after the Wild Card game ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = $offWs.Range("C2").Value() + 26
$offWs.Range("D2").Value = $offWs.Range("D2").Value() + 3
$offWs.Range("E2").Value = $offWs.Range("E2").Value() + 1
$offWs.Range("F2").Value = $offWs.Range("F2").Value() + 9
$offWs.Range("G2").Value = $offWs.Range("G2").Value() + 8
$offWs.Range("J2").Value = $offWs.Range("J2").Value() + 9
$offWs.Range("L2").Value = $offWs.Range("L2").Value() + 71
$offWs.Range("M2").Value = $offWs.Range("M2").Value() + 46
$offWs.Range("O2").Value = $offWs.Range("O2").Value() + 1
$offWs.Range("Q2").Value = $offWs.Range("Q2").Value() + 131
$offWs.Range("B3").Value = $offWs.Range("B3").Value() + 2
$offWs.Range("C3").Value = $offWs.Range("C3").Value() + 31
$offWs.Range("E3").Value = $offWs.Range("E3").Value() + 2
$offWs.Range("F3").Value = $offWs.Range("F3").Value() + 17
$offWs.Range("G3").Value = $offWs.Range("G3").Value() + 4
$offWs.Range("H3").Value = $offWs.Range("H3").Value() + 2
$offWs.Range("I3").Value = $offWs.Range("I3").Value() + 5
$offWs.Range("J3").Value = $offWs.Range("J3").Value() + 7

# --- DEF: season totals after the Wild Card game ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("B2").Value = $defWs.Range("B2").Value() + 1
$defWs.Range("C2").Value = $defWs.Range("C2").Value() + 19
$defWs.Range("D2").Value = $defWs.Range("D2").Value() + 2
$defWs.Range("E2").Value = $defWs.Range("E2").Value() + 2
$defWs.Range("F2").Value = $defWs.Range("F2").Value() + 4
$defWs.Range("G2").Value = $defWs.Range("G2").Value() + 4
$defWs.Range("I2").Value = $defWs.Range("I2").Value() + 2
$defWs.Range("J2").Value = $defWs.Range("J2").Value() + 3
$defWs.Range("L2").Value = $defWs.Range("L2").Value() + 57
$defWs.Range("M2").Value = $defWs.Range("M2").Value() + 31
$defWs.Range("O2").Value = $defWs.Range("O2").Value() + 5
$defWs.Range("P2").Value = $defWs.Range("P2").Value() + 5
$defWs.Range("Q2").Value = $defWs.Range("Q2").Value() + 111
$defWs.Range("B3").Value = $defWs.Range("B3").Value() + 1
$defWs.Range("C3").Value = $defWs.Range("C3").Value() + 18
$defWs.Range("E3").Value = $defWs.Range("E3").Value() + 4
$defWs.Range("F3").Value = $defWs.Range("F3").Value() + 16
$defWs.Range("G3").Value = $defWs.Range("G3").Value() + 3
$defWs.Range("H3").Value = $defWs.Range("H3").Value() + 6
$defWs.Range("I3").Value = $defWs.Range("I3").Value() + 11
$defWs.Range("J3").Value = $defWs.Range("J3").Value() + 6
$defWs.Range("N3").Value = $defWs.Range("N3").Value() + 13

# --- ST: append the Wild Card game log to the running return-yardage sequences ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B4").Value = $stWs.Range("B4").Value() + ' 46 67 60 62 66 64'
$stWs.Range("B5").Value = $stWs.Range("B5").Value() + ' 7 15 16 21 25 26'
$stWs.Range("B6").Value = $stWs.Range("B6").Value() + ' 15 5 15 0'
$stWs.Range("D3").Value = $stWs.Range("D3").Value() + ' 42 28 21 22 49 53 13 43'
$stWs.Range("D4").Value = $stWs.Range("D4").Value() + ' 0 0 0 0 0 3 7 0'
$stWs.Range("D5").Value = $stWs.Range("D5").Value() + ' 0 5 3 17 0 0 9 0 0 0 2 52'

# --- ST: season totals after the Wild Card game ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = $stWs.Range("B2").Value() + 12
$stWs.Range("D2").Value = $stWs.Range("D2").Value() + 8
$stWs.Range("F2").Value = $stWs.Range("F2").Value() + 10
$stWs.Range("G2").Value = $stWs.Range("G2").Value() + 8
$stWs.Range("J2").Value = $stWs.Range("J2").Value() + 1
$stWs.Range("K2").Value = $stWs.Range("K2").Value() + 1
$stWs.Range("L2").Value = $stWs.Range("L2").Value() + 1
$stWs.Range("M2").Value = $stWs.Range("M2").Value() + 1
$stWs.Range("B3").Value = $stWs.Range("B3").Value() + 6

# --- TURNS: season totals after the Wild Card game ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C2").Value = $turnsWs.Range("C2").Value() + 2
$turnsWs.Range("D2").Value = $turnsWs.Range("D2").Value() + 1
$turnsWs.Range("E2").Value = $turnsWs.Range("E2").Value() + 1
$turnsWs.Range("D3").Value = $turnsWs.Range("D3").Value() + 1
$turnsWs.Range("E3").Value = $turnsWs.Range("E3").Value() + 2

# --- PEN: season totals after the Wild Card game ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B3").Value = $penWs.Range("B3").Value() + 1
$penWs.Range("B4").Value = $penWs.Range("B4").Value() + 1
$penWs.Range("D4").Value = $penWs.Range("D4").Value() + 1

